$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 6.636579666666667
$ws.Cells.Item(2, 8).Value = 19.909739
$ws.Cells.Item(2, 9).Value = 0.1201574291771603
$ws.Cells.Item(2, 10).Value = 0.1201574291771603
$ws.Cells.Item(2, 13).Value = 9.467112666666667
$ws.Cells.Item(2, 14).Value = 28.401338
$ws.Cells.Item(2, 15).Value = 0.2314516669582087
$ws.Cells.Item(2, 16).Value = 0.2314516669582087
$ws.Cells.Item(2, 17).Value = 62.82924742564245
$ws.Cells.Item(2, 18).Value = 565.4632268307821
$ws.Cells.Item(2, 19).Value = 0.02781063728046666
$ws.Cells.Item(2, 20).Value = 0.02781063728046666
$ws.Cells.Item(3, 7).Value = 6.636579666666667
$ws.Cells.Item(3, 8).Value = 19.909739
$ws.Cells.Item(3, 9).Value = 0.1201574291771603
$ws.Cells.Item(3, 10).Value = 0.1201574291771603
$ws.Cells.Item(3, 15).Value = 0.03494502079849753
$ws.Cells.Item(3, 16).Value = 0.03494502079849753
$ws.Cells.Item(3, 17).Value = 9.486081422085668
$ws.Cells.Item(3, 18).Value = 85.37473279877101
$ws.Cells.Item(3, 19).Value = 0.00419890386168986
$ws.Cells.Item(3, 20).Value = 0.004198903861689861
$ws.Cells.Item(4, 7).Value = 6.636579666666667
$ws.Cells.Item(4, 8).Value = 19.909739
$ws.Cells.Item(4, 9).Value = 0.1201574291771603
$ws.Cells.Item(4, 10).Value = 0.1201574291771603
$ws.Cells.Item(4, 13).Value = 10.196198
$ws.Cells.Item(4, 14).Value = 30.588594
$ws.Cells.Item(4, 15).Value = 0.2492763218130026
$ws.Cells.Item(4, 16).Value = 0.2492763218130026
$ws.Cells.Item(4, 17).Value = 67.66788032410734
$ws.Cells.Item(4, 18).Value = 609.010922916966
$ws.Cells.Item(4, 19).Value = 0.02995240198378889
$ws.Cells.Item(4, 20).Value = 0.02995240198378889
$ws.Cells.Item(5, 7).Value = 6.636579666666667
$ws.Cells.Item(5, 8).Value = 19.909739
$ws.Cells.Item(5, 9).Value = 0.1201574291771603
$ws.Cells.Item(5, 10).Value = 0.1201574291771603
$ws.Cells.Item(5, 13).Value = 19.81052133333333
$ws.Cells.Item(5, 14).Value = 59.431564
$ws.Cells.Item(5, 15).Value = 0.4843269904302911
$ws.Cells.Item(5, 16).Value = 0.4843269904302911
$ws.Cells.Item(5, 17).Value = 131.4741030668662
$ws.Cells.Item(5, 18).Value = 1183.266927601796
$ws.Cells.Item(5, 19).Value = 0.0581954860512149
$ws.Cells.Item(5, 20).Value = 0.05819548605121491
$ws.Cells.Item(6, 9).Value = 0.00477103065019021
$ws.Cells.Item(6, 10).Value = 0.00477103065019021
$ws.Cells.Item(6, 13).Value = 9.467112666666667
$ws.Cells.Item(6, 14).Value = 28.401338
$ws.Cells.Item(6, 15).Value = 0.2314516669582087
$ws.Cells.Item(6, 16).Value = 0.2314516669582087
$ws.Cells.Item(6, 17).Value = 2.494729350060889
$ws.Cells.Item(6, 18).Value = 22.452564150548
$ws.Cells.Item(6, 19).Value = 0.00110426299709523
$ws.Cells.Item(6, 20).Value = 0.00110426299709523
$ws.Cells.Item(7, 9).Value = 0.00477103065019021
$ws.Cells.Item(7, 10).Value = 0.00477103065019021
$ws.Cells.Item(7, 15).Value = 0.03494502079849753
$ws.Cells.Item(7, 16).Value = 0.03494502079849753
$ws.Cells.Item(7, 19).Value = 0.000166723765301166
$ws.Cells.Item(7, 20).Value = 0.000166723765301166
$ws.Cells.Item(8, 9).Value = 0.00477103065019021
$ws.Cells.Item(8, 10).Value = 0.00477103065019021
$ws.Cells.Item(8, 13).Value = 10.196198
$ws.Cells.Item(8, 14).Value = 30.588594
$ws.Cells.Item(8, 15).Value = 0.2492763218130026
$ws.Cells.Item(8, 16).Value = 0.2492763218130026
$ws.Cells.Item(8, 17).Value = 2.686854514702667
$ws.Cells.Item(8, 18).Value = 24.181690632324
$ws.Cells.Item(8, 19).Value = 0.001189304971736514
$ws.Cells.Item(8, 20).Value = 0.001189304971736514
$ws.Cells.Item(9, 9).Value = 0.00477103065019021
$ws.Cells.Item(9, 10).Value = 0.00477103065019021
$ws.Cells.Item(9, 13).Value = 19.81052133333333
$ws.Cells.Item(9, 14).Value = 59.431564
$ws.Cells.Item(9, 15).Value = 0.4843269904302911
$ws.Cells.Item(9, 16).Value = 0.4843269904302911
$ws.Cells.Item(9, 17).Value = 5.220376132660444
$ws.Cells.Item(9, 18).Value = 46.983385193944
$ws.Cells.Item(9, 19).Value = 0.0023107389160573
$ws.Cells.Item(9, 20).Value = 0.0023107389160573
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.2809586666666666
$ws.Cells.Item(10, 8).Value = 0.842876
$ws.Cells.Item(10, 9).Value = 0.005086847862502274
$ws.Cells.Item(10, 10).Value = 0.005086847862502274
$ws.Cells.Item(10, 13).Value = 9.467112666666667
$ws.Cells.Item(10, 14).Value = 28.401338
$ws.Cells.Item(10, 15).Value = 0.2314516669582087
$ws.Cells.Item(10, 16).Value = 0.2314516669582087
$ws.Cells.Item(10, 17).Value = 2.659867352009778
$ws.Cells.Item(10, 18).Value = 23.938806168088
$ws.Cells.Item(10, 19).Value = 0.001177359417338952
$ws.Cells.Item(10, 20).Value = 0.001177359417338952
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.2809586666666666
$ws.Cells.Item(11, 8).Value = 0.842876
$ws.Cells.Item(11, 9).Value = 0.005086847862502274
$ws.Cells.Item(11, 10).Value = 0.005086847862502274
$ws.Cells.Item(11, 15).Value = 0.03494502079849753
$ws.Cells.Item(11, 16).Value = 0.03494502079849753
$ws.Cells.Item(11, 17).Value = 0.4015919226626667
$ws.Cells.Item(11, 18).Value = 3.614327303964
$ws.Cells.Item(11, 19).Value = 0.0001777600043539346
$ws.Cells.Item(11, 20).Value = 0.0001777600043539346
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.2809586666666666
$ws.Cells.Item(12, 8).Value = 0.842876
$ws.Cells.Item(12, 9).Value = 0.005086847862502274
$ws.Cells.Item(12, 10).Value = 0.005086847862502274
$ws.Cells.Item(12, 13).Value = 10.196198
$ws.Cells.Item(12, 14).Value = 30.588594
$ws.Cells.Item(12, 15).Value = 0.2492763218130026
$ws.Cells.Item(12, 16).Value = 0.2492763218130026
$ws.Cells.Item(12, 17).Value = 2.864710195149333
$ws.Cells.Item(12, 18).Value = 25.782391756344
$ws.Cells.Item(12, 19).Value = 0.001268030724786901
$ws.Cells.Item(12, 20).Value = 0.001268030724786901
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.2809586666666666
$ws.Cells.Item(13, 8).Value = 0.842876
$ws.Cells.Item(13, 9).Value = 0.005086847862502274
$ws.Cells.Item(13, 10).Value = 0.005086847862502274
$ws.Cells.Item(13, 13).Value = 19.81052133333333
$ws.Cells.Item(13, 14).Value = 59.431564
$ws.Cells.Item(13, 15).Value = 0.4843269904302911
$ws.Cells.Item(13, 16).Value = 0.4843269904302911
$ws.Cells.Item(13, 17).Value = 5.565937659784888
$ws.Cells.Item(13, 18).Value = 50.093438938064
$ws.Cells.Item(13, 19).Value = 0.002463697716022486
$ws.Cells.Item(13, 20).Value = 0.002463697716022486
$ws.Cells.Item(14, 7).Value = 48.051317
$ws.Cells.Item(14, 8).Value = 144.153951
$ws.Cells.Item(14, 9).Value = 0.8699846923101473
$ws.Cells.Item(14, 10).Value = 0.8699846923101472
$ws.Cells.Item(14, 13).Value = 9.467112666666667
$ws.Cells.Item(14, 14).Value = 28.401338
$ws.Cells.Item(14, 15).Value = 0.2314516669582087
$ws.Cells.Item(14, 16).Value = 0.2314516669582087
$ws.Cells.Item(14, 17).Value = 454.9072318207154
$ws.Cells.Item(14, 18).Value = 4094.165086386438
$ws.Cells.Item(14, 19).Value = 0.2013594072633079
$ws.Cells.Item(14, 20).Value = 0.2013594072633079
$ws.Cells.Item(15, 7).Value = 48.051317
$ws.Cells.Item(15, 8).Value = 144.153951
$ws.Cells.Item(15, 9).Value = 0.8699846923101473
$ws.Cells.Item(15, 10).Value = 0.8699846923101472
$ws.Cells.Item(15, 15).Value = 0.03494502079849753
$ws.Cells.Item(15, 16).Value = 0.03494502079849753
$ws.Cells.Item(15, 17).Value = 68.68277462107102
$ws.Cells.Item(15, 18).Value = 618.1449715896391
$ws.Cells.Item(15, 19).Value = 0.03040163316715257
$ws.Cells.Item(15, 20).Value = 0.03040163316715256
$ws.Cells.Item(16, 7).Value = 48.051317
$ws.Cells.Item(16, 8).Value = 144.153951
$ws.Cells.Item(16, 9).Value = 0.8699846923101473
$ws.Cells.Item(16, 10).Value = 0.8699846923101472
$ws.Cells.Item(16, 13).Value = 10.196198
$ws.Cells.Item(16, 14).Value = 30.588594
$ws.Cells.Item(16, 15).Value = 0.2492763218130026
$ws.Cells.Item(16, 16).Value = 0.2492763218130026
$ws.Cells.Item(16, 17).Value = 489.9407422927661
$ws.Cells.Item(16, 18).Value = 4409.466680634894
$ws.Cells.Item(16, 19).Value = 0.2168665841326904
$ws.Cells.Item(16, 20).Value = 0.2168665841326903
$ws.Cells.Item(17, 7).Value = 48.051317
$ws.Cells.Item(17, 8).Value = 144.153951
$ws.Cells.Item(17, 9).Value = 0.8699846923101473
$ws.Cells.Item(17, 10).Value = 0.8699846923101472
$ws.Cells.Item(17, 13).Value = 19.81052133333333
$ws.Cells.Item(17, 14).Value = 59.431564
$ws.Cells.Item(17, 15).Value = 0.4843269904302911
$ws.Cells.Item(17, 16).Value = 0.4843269904302911
$ws.Cells.Item(17, 17).Value = 951.9216405232628
$ws.Cells.Item(17, 18).Value = 8567.294764709364
$ws.Cells.Item(17, 19).Value = 0.4213570677469965
$ws.Cells.Item(17, 20).Value = 0.4213570677469964
